$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep their original text formatting
# (values such as "218.30", "7.850", "0.05150" would otherwise be auto-converted
# to numbers by Excel, stripping significant trailing zeros / dot grouping).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.054.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.652.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5284"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.52%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2621"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06314"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.713.13"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.81%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.488"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5469"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.067.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.570"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.12"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.48%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "139.78"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1245"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.276"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05945"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.275"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.251"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.541"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.415"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9431"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.756"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5642"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01611"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.888"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8479"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.94"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007.70"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.790.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.90"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.487"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4287"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.37%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.850"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.97%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05150"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.57%  "
